# Update LDLC prices history
#
# A new price-history column is inserted right before the existing "nom"
# column (BR). This shifts the old BR ("nom") -> BS and the old BS
# ("url_produit") -> BT, and the new BR column is filled with:
#   - row 1 (header): the new snapshot timestamp
#   - rows 2-80: a copy of the latest existing price snapshot (column BQ)
#   - rows 81-206: left blank, matching the already-blank BQ cells on
#     those rows (products with no price recorded in this snapshot)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at BR; everything from BR onward shifts one column
# to the right (BR -> BS, BS -> BT), picking up the formatting of the
# inserted-before column as Excel normally does.
$ws.Columns("BR:BR").Insert()

# New header cell for the freshly inserted price-snapshot column.
$ws.Range("BR1").Value2 = "2026-01-30 20:16:43"

# Rows 2-80 carry a numeric price in BQ; duplicate it into the new BR
# column for this snapshot (prices unchanged since the previous poll).
$ws.Range("BR2:BR80").Value2 = $ws.Range("BQ2:BQ80").Value2
